$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows above the old row 5 ("state"), pushing everything
# else down by two rows.
$ws.Rows("5:6").Insert() | Out-Null

# New row 5: occupation / F52
$ws.Range("A5").Value = "occupation"
# New row 6: industry / F53
$ws.Range("A6").Value = "industry"

# Fill in the F53 (industry) text cells before the F52 (occupation) ones so
# the shared-string table ends up ordered the same way the original author's
# edit produced it (occupation, industry, F53, F52).
$ws.Range("C6").Value = "F53"
$ws.Range("E6").Value = "F53"

$ws.Range("C5").Value = "F52"
$ws.Range("E5").Value = "F52"

$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("D6").Value = 1

# Leave the selection where the author left it after the edit.
$ws.Range("C6").Select() | Out-Null
